$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.971499999999988
$ws.Range("D4").Value = -7.288199999999997
$ws.Range("A9").Value = -20.49569999999997
$ws.Range("D10").Value = -7.679399999999995
$ws.Range("E12").Value = 12.22839999999999
$ws.Range("E17").Value = 13.57350000000001
$ws.Range("A18").Value = -22.85640000000001
$ws.Range("A20").Value = -22.11570000000002
$ws.Range("C21").Value = -13.27090000000001
